$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clinical_assessment_categories")

# New column E: clinical_assessment_category_order_all
$ws.Range("E1").Value = "clinical_assessment_category_order_all"
# Give the new header cell the same distinctive (black Calibri) font used by
# the other header cells in row 1 (C1/D1).
$ws.Range("E1").Font.Name = "Calibri"

$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 4

# Column E should share the same width as column D.
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Match the view state recorded in the saved workbook: this sheet becomes
# the active tab with E6 selected.
$ws.Activate()
$ws.Range("E6").Select()
